$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of column A with column E, and column B with column F,
# for each of the 25 data rows (A1:F25).
for ($r = 1; $r -le 25; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)
    $cellE = $ws.Cells.Item($r, 5)
    $cellF = $ws.Cells.Item($r, 6)

    $valA = $cellA.Value2
    $valB = $cellB.Value2
    $valE = $cellE.Value2
    $valF = $cellF.Value2

    # Only write back when the source cell actually holds a value, so that
    # already-empty cells are left untouched (writing an empty value to an
    # empty cell would otherwise delete the cell node entirely).
    if ($valE -ne $null -and $valE -ne "") {
        $cellA.Value2 = $valE
    }
    if ($valA -ne $null -and $valA -ne "") {
        $cellE.Value2 = $valA
    }

    if ($valF -ne $null -and $valF -ne "") {
        $cellB.Value2 = $valF
    }
    if ($valB -ne $null -and $valB -ne "") {
        $cellF.Value2 = $valB
    }
}
